$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Status" (column D) values: mark several topics as "done" and one as "revisit" ---
$ws.Range("D8").Value = "done"
$ws.Range("D10").Value = "revisit"
$ws.Range("D15").Value = "done"

# --- Row 16: fill in the "learnings" note + status ---
$ws.Range("C16").Value = "variable declaration are made undefined, function scope"
$ws.Range("D16").Value = "done"

# Row 16's D cell should pick up the same formatting as the rest of that row
# (copy format only, matching how the rest of column C/D is styled).
$ws.Range("C12").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null

# --- Row 17: fill in the "learnings" note + status ---
$ws.Range("C17").Value = "function defined as variable will always not have hoisting feature and will remain undefined "
$ws.Range("D17").Value = "done"

# --- Update the view: scroll position, zoom, and active selection ---
$ws.Range("A12").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 153
$ws.Range("B19").Select()
